# Applies the "update data files, remove redundant old versions" edit:
# adds four new max/min capacity-potential columns (BE:BL) with per-year
# (2020/2030/2040/2050-ish "b1..b4") bucket values for each conversion
# technology row, plus their header labels, units row and an "e for
# electricity" annotation, and removes the now-redundant placeholder
# blank cells that used to live in BD:BH on the remarks/sources rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New annotation text + unit labels.
#    Written first so the shared-string table allocates indices in the
#    same order the source workbook ended up with (120, then 121-128 for
#    the header row, then 129 for the "MW" unit labels).
# ---------------------------------------------------------------------
$ws.Range("B21").Value = "e for electricity"

# ---------------------------------------------------------------------
# 2) New header row (row 1) labels for columns BE:BL.
# ---------------------------------------------------------------------
$ws.Range("BE1").Value = "maxCapacityPotCTb1"
$ws.Range("BF1").Value = "maxCapacityPotCTb2"
$ws.Range("BG1").Value = "maxCapacityPotCTb3"
$ws.Range("BH1").Value = "maxCapacityPotCTb4"
$ws.Range("BI1").Value = "minCapacityPotCTb1"
$ws.Range("BJ1").Value = "minCapacityPotCTb2"
$ws.Range("BK1").Value = "minCapacityPotCTb3"
$ws.Range("BL1").Value = "minCapacityPotCTb4"

# ---------------------------------------------------------------------
# 3) Units row labels ("MW") under the Remarks block.
# ---------------------------------------------------------------------
$ws.Range("BE21").Value = "MW"
$ws.Range("BI21").Value = "MW"

# ---------------------------------------------------------------------
# 4) Per-technology data, rows 2-19 (BE:BH = max capacity potential
#    buckets, BI:BL = min capacity potential buckets).
# ---------------------------------------------------------------------
$maxVals = @{
  2  = 176979.12
  3  = 176979.12
  4  = 10000
  5  = 10000
  6  = 10000
  7  = 10000
  8  = 50000
  9  = 50000
  10 = 50000
  11 = 50000
  12 = 50000
  13 = 50000
  14 = 50000
  15 = 50000
  16 = 50000
  17 = 50000
  18 = 10000
  19 = 10000
}
$maxValsF = @{
  2  = 162660.01999999999
  3  = 162660.01999999999
}
$maxValsG = @{
  2  = 112895.14
  3  = 112895.14
}
$maxValsH = @{
  2  = 145846.91
  3  = 145846.91
}

foreach ($r in 2..19) {
  $be = $maxVals[$r]
  $bf = if ($maxValsF.ContainsKey($r)) { $maxValsF[$r] } else { $maxVals[$r] }
  $bg = if ($maxValsG.ContainsKey($r)) { $maxValsG[$r] } else { $maxVals[$r] }
  $bh = if ($maxValsH.ContainsKey($r)) { $maxValsH[$r] } else { $maxVals[$r] }

  $ws.Range("BE$r").Value = $be
  $ws.Range("BF$r").Value = $bf
  $ws.Range("BG$r").Value = $bg
  $ws.Range("BH$r").Value = $bh

  $ws.Range("BI$r").Value = 0
  $ws.Range("BJ$r").Value = 0
  $ws.Range("BK$r").Value = 0
  $ws.Range("BL$r").Value = 0

  # Rows 2 & 3 keep the default (General) number format on BE:BH, every
  # other row (and BI:BL on every row) uses the workbook's "0.0" style.
  if ($r -ne 2 -and $r -ne 3) {
    $ws.Range("BE$r`:BH$r").NumberFormat = "0.0"
  }
  $ws.Range("BI$r`:BL$r").NumberFormat = "0.0"
}

# ---------------------------------------------------------------------
# 5) Clean up the now-redundant blank placeholder cells (old style-only
#    filler cells) that used to span BD:BH on rows 22-34 (and BE:BH on
#    rows 21-22, 30-34) -- BI onward stays untouched.
# ---------------------------------------------------------------------
$ws.Range("BF21:BH21").Clear()
$ws.Range("BE22:BH22").Clear()
foreach ($r in 23..29) {
  $ws.Range("BD$r`:BH$r").Clear()
}
foreach ($r in 30..34) {
  $ws.Range("BE$r`:BH$r").Clear()
}

# ---------------------------------------------------------------------
# 6) Selection / scroll position, best effort.
# ---------------------------------------------------------------------
$ws.Range("BJ24").Select()
